# Updates cryptocurrency price/volume figures per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.530.94'
$ws.Range("E2").Value = '  +2.01%  '
$ws.Range("D3").Value = '1.909.84'
$ws.Range("E3").Value = '  +5.32%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.93'
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5042'
$ws.Range("E7").Value = '  +1.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3955'
$ws.Range("E8").Value = '  +1.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09694'
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.163'
$ws.Range("E10").Value = '  +5.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.53'
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.546'
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.11'
$ws.Range("E13").Value = '  +2.91%  '
$ws.Range("D14").Value = '1.920.39'
$ws.Range("E14").Value = '  +5.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.574'
$ws.Range("E15").Value = '  +3.66%  '
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001145'
$ws.Range("E17").Value = '  +1.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.87'
$ws.Range("E18").Value = '  +1.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06652'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("E20").Value = '  +5.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.291'
$ws.Range("E22").Value = '  +6.24%  '
$ws.Range("D23").Value = '28.585.93'
$ws.Range("E23").Value = '  +2.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.45'
$ws.Range("E24").Value = '  +3.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.283'
$ws.Range("E25").Value = '  +1.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.760'
$ws.Range("E26").Value = '  +15.29%  '
$ws.Range("D27").Value = '2.135.04'
$ws.Range("E27").Value = '  +5.66%  '
$ws.Range("E28").Value = '  +3.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '159.65'
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.86'
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.108'
$ws.Range("E31").Value = '  +6.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1072'
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.710'
$ws.Range("E33").Value = '  +2.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.641'
$ws.Range("E34").Value = '  +0.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.936'
$ws.Range("E35").Value = '  +10.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06776'
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02450'
$ws.Range("E37").Value = '  +4.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2234'
$ws.Range("E38").Value = '  +4.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.091'
$ws.Range("E39").Value = '  +2.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.64'
$ws.Range("E40").Value = '  +3.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6439'
$ws.Range("E41").Value = '  +4.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.192'
$ws.Range("E42").Value = '  +3.94%  '
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.68'
$ws.Range("E44").Value = '  +4.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6107'
$ws.Range("E45").Value = '  +3.96%  '
$ws.Range("E46").Value = '  +0.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.668'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.045'
$ws.Range("E48").Value = '  +5.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.22'
$ws.Range("E49").Value = '  +1.98%  '
$ws.Range("E50").Value = '  +2.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.33'
$ws.Range("E51").Value = '  +6.12%  '
